$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Aula 60 - Validacao de data com Spring Validator (Importante para Regras de
# Negocio). Adds three new rows (65-67) to the notes table, mirroring the
# existing layout used for aula 59 (row 64) and earlier lessons.
# ---------------------------------------------------------------------------

# Copy the formatting (font/fill/alignment) of the previous lesson's row down
# onto the three new rows so the B/C columns reuse the "s=5" style and the
# D/E columns reuse the "s=1" (wrap text) style already present in the sheet.
$ws.Range("B64:E64").Copy() | Out-Null
$ws.Range("B65:E65").PasteSpecial(-4122) | Out-Null
$ws.Range("B66:E66").PasteSpecial(-4122) | Out-Null
$ws.Range("B67:E67").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$tituloAula60 = "`r`n60. Validação de data com Spring Validator"

$notaInitBinder = "8:21`r`nanotação @InitBinder : está anotação diz para a aplicação que este método sera o primeiro metodo da classe que vai ser executado ao ser chamado FuncionarioController. Desta forma, ao ser executado, o SPring MVC vai até a classe FuncionarioValidator fazer a validação antes de liberar o acesso a requisição pra os metodos salvar e editar"

$notaImportante = "10:22`r`nIMPORTANTE - REGRAS DE NEGÓCIO E VALIDAÇÃO: no metodo validate() pode conter qualquer tipo de validação de campos HTML, perfeito para validar as regras de negócio da aplicação."

$notaMessages = "7:11`r`na criação do arquivo messages.properties é devido a especificação do Spring Validator`r`ne o arquivo ValidationMessages.Properties é um arquivo próprio`r`nda especificação Bean Validation, que automaticamente já procura este arquivo no classpath`r`no Spring Validation foi criado antes do Bean validation"

# Cell values are written in the same order the author typed/pasted them so
# that new entries land on the shared-string table in the matching order
# (112=title, 113=InitBinder note, 114=IMPORTANTE note, 115=messages note).
$ws.Range("B65").Value = 60
$ws.Range("C65").Value = "11. Validação Back-End"
$ws.Range("D65").Value = $tituloAula60

$ws.Range("B66").Value = 60
$ws.Range("C66").Value = "11. Validação Back-End"
$ws.Range("D66").Value = $tituloAula60
$ws.Range("E66").Value = $notaInitBinder

$ws.Range("B67").Value = 60
$ws.Range("C67").Value = "11. Validação Back-End"
$ws.Range("D67").Value = $tituloAula60
$ws.Range("E67").Value = $notaImportante
$ws.Range("E67").WrapText = $true
$ws.Range("E67").Font.Bold = $true
$ws.Range("E67").Interior.Color = 49407

$ws.Range("E65").Value = $notaMessages

# Row heights, matching what Excel computed for the wrapped text in each row
$ws.Rows.Item(65).RowHeight = 120
$ws.Rows.Item(66).RowHeight = 45
$ws.Rows.Item(67).RowHeight = 45

# Scroll/selection state left by the author after typing the new notes
$ws.Range("E66").Select() | Out-Null
